# Insert a new data row at row 108, pushing all existing rows 108..170 down to 109..171.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("108").Insert()

# Populate the new row 108 with the new record (weekly update for Chirimoya prices)
$ws.Range("A108").Value = 10
$ws.Range("B108").Value = "Vega Modelo de Temuco"
$ws.Range("C108").Value = "La Araucanía"
$ws.Range("D108").Value = 44879
$ws.Range("E108").Value = 9
$ws.Range("F108").Value = "Fruta"
$ws.Range("G108").Value = 100107
$ws.Range("H108").Value = "Otros"
$ws.Range("I108").Value = 100107002
$ws.Range("J108").Value = "Chirimoya"
$ws.Range("K108").Value = "Cultivar IV Región"
$ws.Range("L108").Value = "Primera"
$ws.Range("M108").Value = 50
$ws.Range("N108").Value = 2800
$ws.Range("O108").Value = 2800
$ws.Range("P108").Value = 2800
$ws.Range("Q108").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R108").Value = "Provincia del Elquí"
$ws.Range("S108").Value = 2800
$ws.Range("T108").Value = 1
